$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for data rows 2 through 32
# from serial 45609 (2024-11-13) to serial 45610 (2024-11-14)
for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 3).Value = 45610
}

# Rows 30 and 31 had their "Beteckning" (A) and "Area (ha)" (G) values swapped
$ws.Range("A30").Value = "A 46082-2024"
$ws.Range("G30").Value = 1

$ws.Range("A31").Value = "A 46085-2024"
$ws.Range("G31").Value = 1.9
